# Update the emergency-room intake form with the new patient's data
# (commit: "Version 2.0.1 solucionado error espera de base de datos").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient name / clinical file number (row 6) ---
$ws.Range("A6").Value = "CALABERA   HERNANDEZ  JIMMY  EMANUEL"
$ws.Range("G6").Value = "/201761945"

# --- Birth date / age / place of birth / sex (row 9) ---
# Leading apostrophe forces text so Excel doesn't coerce the
# ISO-looking date / plain number into a date serial or numeric value.
$ws.Range("A9").Value = "'2004-06-17"
$ws.Range("D9").Value = "'13"
$ws.Range("E9").Value = "GUATEMALA"
$ws.Range("G9").Value = "MASCULINO"

# --- Occupation / Nationality / ID document (row 11) ---
$ws.Range("C11").Value = "ESTUDIANTE "
$ws.Range("E11").Value = "GUATEMALA "
$ws.Range("G11").Value = "'"

# --- Emergency contact (row 13) ---
$ws.Range("A13").Value = "JIMMY CALABERA "
$ws.Range("D13").Value = "PAPA"
$ws.Range("E13").Value = "COL. ILUSIONES SECCION K LOTE 16"
$ws.Range("G13").Value = "42 85 12 27"

# --- Time of medical assistance (row 14) ---
$ws.Range("D14").Value = "Hora: 12:18:7"
